$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 3.25
$ws.Range("L2").Value = 1.48
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 1.39
$ws.Range("P2").Value = 1.79
$ws.Range("Q2").Value = 2.16
$ws.Range("R2").Value = 1.3
$ws.Range("S2").Value = 3.9
$ws.Range("T2").Value = 1.85
$ws.Range("U2").Value = 2.08
$ws.Range("W2").Value = 1.67
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 12.5
$ws.Range("AA2").Value = 170
$ws.Range("AB2").Value = 9.6
$ws.Range("AD2").Value = 15
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 18.5
$ws.Range("AI2").Value = 60
$ws.Range("AK2").Value = 27
$ws.Range("AL2").Value = 130
$ws.Range("AM2").Value = 120
$ws.Range("AN2").Value = 30
$ws.Range("AO2").Value = 50

# Row 3
$ws.Range("F3").Value = 1.72
$ws.Range("G3").Value = 1.79
$ws.Range("I3").Value = 6.8
$ws.Range("J3").Value = 3.6
$ws.Range("L3").Value = 1.57
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 3
$ws.Range("P3").Value = 1.65
$ws.Range("Q3").Value = 2.36
$ws.Range("R3").Value = 1.24
$ws.Range("S3").Value = 4.7
$ws.Range("T3").Value = 2.16
$ws.Range("U3").Value = 1.72
$ws.Range("Y3").Value = 980
$ws.Range("Z3").Value = 1000
$ws.Range("AH3").Value = 980
$ws.Range("AL3").Value = 300

# Row 4
$ws.Range("F4").Value = 2.92
$ws.Range("G4").Value = 2.96
$ws.Range("H4").Value = 2.88
$ws.Range("L4").Value = 1.52
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 3.05
$ws.Range("S4").Value = 4.7
$ws.Range("T4").Value = 1.92
$ws.Range("U4").Value = 1.94
$ws.Range("W4").Value = 1.51
$ws.Range("X4").Value = 10.5
$ws.Range("Y4").Value = 9.6
$ws.Range("Z4").Value = 17.5
$ws.Range("AA4").Value = 46
$ws.Range("AB4").Value = 9.6
$ws.Range("AC4").Value = 7.4
$ws.Range("AD4").Value = 13
$ws.Range("AE4").Value = 38
$ws.Range("AF4").Value = 19.5
$ws.Range("AG4").Value = 13.5
$ws.Range("AJ4").Value = 50
$ws.Range("AK4").Value = 40
$ws.Range("AL4").Value = 60
$ws.Range("AM4").Value = 140
$ws.Range("AN4").Value = 55
$ws.Range("AO4").Value = 40

# Row 5
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 6.4
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4.2
$ws.Range("L5").Value = 1.48
$ws.Range("N5").Value = 3.2
$ws.Range("P5").Value = 1.74
$ws.Range("Q5").Value = 2.24
$ws.Range("R5").Value = 1.27
$ws.Range("W5").Value = 2.42
$ws.Range("AF5").Value = 10
$ws.Range("AN5").Value = 16

# Row 6
$ws.Range("F6").Value = 1.94
$ws.Range("G6").Value = 1.98
$ws.Range("H6").Value = 4.7
$ws.Range("I6").Value = 5.1
$ws.Range("J6").Value = 3.25
$ws.Range("K6").Value = 3.5
$ws.Range("L6").Value = 1.49
$ws.Range("M6").Value = 1.09
$ws.Range("N6").Value = 3.25
$ws.Range("O6").Value = 1.41
$ws.Range("P6").Value = 1.76
$ws.Range("Q6").Value = 2.22
$ws.Range("R6").Value = 1.29
$ws.Range("S6").Value = 4.1
$ws.Range("U6").Value = 1.88
$ws.Range("V6").Value = 1.24
$ws.Range("W6").Value = 2.02
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 15
$ws.Range("AA6").Value = 120
$ws.Range("AB6").Value = 8.800000000000001
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 19.5
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 11.5
$ws.Range("AH6").Value = 22
$ws.Range("AI6").Value = 100
$ws.Range("AJ6").Value = 23
$ws.Range("AK6").Value = 25
$ws.Range("AL6").Value = 46
$ws.Range("AM6").Value = 150
$ws.Range("AO6").Value = 1000

# Row 7
$ws.Range("G7").Value = 1.1
$ws.Range("I7").Value = 1000
$ws.Range("L7").Value = 1.21
$ws.Range("N7").Value = 7.4
$ws.Range("O7").Value = 1.12
$ws.Range("P7").Value = 3.35
$ws.Range("Q7").Value = 1.38
$ws.Range("R7").Value = 1.95
$ws.Range("S7").Value = 1.95
$ws.Range("U7").Value = 1.38
$ws.Range("AB7").Value = 980
$ws.Range("AF7").Value = 9.4
$ws.Range("AJ7").Value = 8.800000000000001

# Row 8
$ws.Range("F8").Value = 1.63
$ws.Range("G8").Value = 1.69
$ws.Range("H8").Value = 6.6
$ws.Range("I8").Value = 7.2
$ws.Range("K8").Value = 4.1
$ws.Range("L8").Value = 1.52
$ws.Range("N8").Value = 2.92
$ws.Range("P8").Value = 1.65
$ws.Range("U8").Value = 1.67
$ws.Range("W8").Value = 2.44
$ws.Range("Y8").Value = 16.5
$ws.Range("AB8").Value = 6.4
$ws.Range("AD8").Value = 27
$ws.Range("AE8").Value = 150
$ws.Range("AF8").Value = 8.4
$ws.Range("AG8").Value = 10.5
$ws.Range("AI8").Value = 140
$ws.Range("AJ8").Value = 16.5
$ws.Range("AN8").Value = 15
$ws.Range("AO8").Value = 1000

# Row 9
$ws.Range("F9").Value = 1.45
$ws.Range("G9").Value = 1.51
$ws.Range("K9").Value = 4.9
$ws.Range("L9").Value = 1.43
$ws.Range("P9").Value = 1.82
$ws.Range("Q9").Value = 2.04
$ws.Range("R9").Value = 1.31
$ws.Range("S9").Value = 3.7
$ws.Range("T9").Value = 2.26
$ws.Range("U9").Value = 1.67
$ws.Range("W9").Value = 2.96
$ws.Range("X9").Value = 16.5
$ws.Range("Y9").Value = 29
$ws.Range("AA9").Value = 460
$ws.Range("AB9").Value = 8
$ws.Range("AC9").Value = 12.5
$ws.Range("AD9").Value = 42
$ws.Range("AE9").Value = 220
$ws.Range("AF9").Value = 9.4
$ws.Range("AG9").Value = 12.5
$ws.Range("AH9").Value = 36
$ws.Range("AI9").Value = 190
$ws.Range("AJ9").Value = 15.5
$ws.Range("AK9").Value = 22
$ws.Range("AL9").Value = 60
$ws.Range("AM9").Value = 260
$ws.Range("AN9").Value = 11.5

# Row 10
$ws.Range("F10").Value = 1.34
$ws.Range("G10").Value = 1.4
$ws.Range("H10").Value = 10
$ws.Range("I10").Value = 12.5
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 5.9
$ws.Range("L10").Value = 1.32
$ws.Range("M10").Value = 1.04
$ws.Range("N10").Value = 1.11
$ws.Range("O10").Value = 1.21
$ws.Range("P10").Value = 1.25
$ws.Range("Q10").Value = 1.54
$ws.Range("R10").Value = 1.46
$ws.Range("S10").Value = 2.6
$ws.Range("T10").Value = 1.03
$ws.Range("U10").Value = 1.03
$ws.Range("W10").Value = 3.35
$ws.Range("X10").Value = 1000
$ws.Range("Y10").Value = 48
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 12.5
$ws.Range("AC10").Value = 1000
$ws.Range("AD10").Value = 980
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 12
$ws.Range("AG10").Value = 15
$ws.Range("AH10").Value = 34
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 15.5
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 55
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000

# Row 11
$ws.Range("N11").Value = 1.89
$ws.Range("O11").Value = 1.22
$ws.Range("P11").Value = 1.89
$ws.Range("Q11").Value = 1.59
$ws.Range("R11").Value = 1.34
$ws.Range("S11").Value = 2.4

# Row 12
$ws.Range("F12").Value = 2.44
$ws.Range("G12").Value = 2.56
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 3.25
$ws.Range("K12").Value = 3.45
$ws.Range("N12").Value = 2.88
$ws.Range("P12").Value = 1.63
$ws.Range("S12").Value = 4.6
$ws.Range("T12").Value = 1.92
$ws.Range("V12").Value = 1.41
$ws.Range("W12").Value = 1.64
$ws.Range("X12").Value = 10
$ws.Range("Y12").Value = 12
$ws.Range("Z12").Value = 21
$ws.Range("AA12").Value = 65
$ws.Range("AB12").Value = 9.800000000000001
$ws.Range("AC12").Value = 8
$ws.Range("AD12").Value = 15
$ws.Range("AF12").Value = 20
$ws.Range("AG12").Value = 13
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 70
$ws.Range("AJ12").Value = 38
$ws.Range("AM12").Value = 150
